# Auto-generated edit script: refresh scraped horarios data (new scrape @ 07:45:49)
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$data1 = @(
    @("03:55:13","04:01","81_EL PELIGRO",6,"LP1912"),
    @("03:55:13","04:34","215_ALUAR",39,"LP1912"),
    @("04:18:52","04:45","215A_EL PATO",27,"LP1912"),
    @("03:55:13","04:53","11_ETCHEVERRY",58,"LP1912"),
    @("04:54:17","04:54","11_ETCHEVERRY",0,"LP1912"),
    @("03:55:13","05:16","17_ROMERO",81,"LP1912"),
    @("04:18:52","05:21","23_HERNANDEZ",63,"LP1912"),
    @("03:55:13","05:22","23_HERNANDEZ",87,"LP1912"),
    @("04:18:52","05:34","215B_EL PATO",76,"LP1912"),
    @("03:55:13","05:35","215B_EL PATO",100,"LP1912"),
    @("03:55:13","05:46","15_ABASTO",111,"LP1912"),
    @("04:18:52","05:53","10_OLMOS",95,"LP1912"),
    @("04:40:48","05:54","10_OLMOS",74,"LP1912"),
    @("04:40:48","06:04","16_SANTA ANA",84,"LP1912"),
    @("04:18:52","06:05","16_SANTA ANA",107,"LP1912"),
    @("04:18:52","06:11","215A_EL PATO",113,"LP1912"),
    @("04:18:52","06:13","225_HARAS DEL SUR",115,"LP1912"),
    @("04:40:48","06:14","225_HARAS DEL SUR",94,"LP1912"),
    @("04:40:48","06:21","26_HERNANDEZ",101,"LP1912"),
    @("04:40:48","06:27","23_HERNANDEZ",107,"LP1912"),
    @("04:40:48","06:29","86_EST CHICA-ESC AGRARIA",109,"LP1912"),
    @("05:49:10","06:30","86_EST CHICA-ESC AGRARIA",41,"LP1912"),
    @("04:40:48","06:31","16_SANTA ANA",111,"LP1912"),
    @("04:54:17","06:44","225_C ROCA-H SUR",110,"LP1912"),
    @("04:54:17","06:46","215C_EL PATO",112,"LP1912"),
    @("05:49:10","06:47","215C_EL PATO",58,"LP1912"),
    @("05:19:24","06:59","14_ABASTO",100,"LP1912"),
    @("05:49:10","07:00","14_ABASTO",71,"LP1912"),
    @("06:14:19","07:01","16_SANTA ANA",47,"LP1912"),
    @("05:19:24","07:05","15_ABASTO",106,"LP1912"),
    @("05:49:10","07:05","23_HERNANDEZ",76,"LP1912"),
    @("05:19:24","07:07","225_GOMEZ",108,"LP1912"),
    @("05:19:24","07:11","215A_EL PATO",112,"LP1912"),
    @("05:19:24","07:15","11_ETCHEVERRY",116,"LP1912"),
    @("05:49:10","07:16","11_ETCHEVERRY",87,"LP1912"),
    @("06:43:12","07:17","16_SANTA ANA",34,"LP1912"),
    @("07:19:37","07:19","10_OLMOS",0,"LP1912"),
    @("05:49:10","07:21","26_HERNANDEZ",92,"LP1912"),
    @("06:14:19","07:23","10_OLMOS",69,"LP1912"),
    @("07:19:37","07:31","16_SANTA ANA",12,"LP1912"),
    @("06:14:19","07:31","11_ETCHEVERRY",77,"LP1912"),
    @("05:49:10","07:32","11_ETCHEVERRY",103,"LP1912"),
    @("05:49:10","07:32","84_COLONIA URQUIZA-ESC 49",103,"LP1912"),
    @("05:49:10","07:32","16_SANTA ANA",103,"LP1912"),
    @("07:19:37","07:34","23_HERNANDEZ",15,"LP1912"),
    @("07:19:37","07:36","27_EL RETIRO",17,"LP1912"),
    @("05:49:10","07:37","27_EL RETIRO",108,"LP1912"),
    @("05:49:10","07:39","10_OLMOS",110,"LP1912"),
    @("07:45:49","07:46","16_SANTA ANA",1,"LP1912"),
    @("06:14:19","07:47","14_ABASTO",93,"LP1912"),
    @("05:49:10","07:48","14_ABASTO",119,"LP1912"),
    @("06:14:19","07:51","215D_EL PATO",97,"LP1912"),
    @("06:43:12","07:52","215D_EL PATO",69,"LP1912"),
    @("07:45:49","07:55","10_OLMOS",10,"LP1912"),
    @("07:19:37","07:59","23_HERNANDEZ",40,"LP1912"),
    @("06:14:19","08:00","23_HERNANDEZ",106,"LP1912"),
    @("07:19:37","08:03","11_ETCHEVERRY",44,"LP1912"),
    @("06:43:12","08:03","23_HERNANDEZ",80,"LP1912"),
    @("06:57:11","08:05","23_HERNANDEZ",68,"LP1912"),
    @("07:19:37","08:10","16_SANTA ANA",51,"LP1912"),
    @("06:14:19","08:12","15_ABASTO",118,"LP1912"),
    @("07:45:49","08:13","10_OLMOS",28,"LP1912"),
    @("06:43:12","08:21","26_HERNANDEZ",98,"LP1912"),
    @("07:19:37","08:22","16_P MOR-SANTA ANA",63,"LP1912"),
    @("06:43:12","08:23","215B_EL PATO",100,"LP1912"),
    @("06:43:12","08:23","16_P MOR-SANTA ANA",100,"LP1912"),
    @("06:43:12","08:27","84_COLONIA URQUIZA-ESC 49",104,"LP1912"),
    @("07:45:49","08:31","23_HERNANDEZ",46,"LP1912"),
    @("06:57:11","08:42","81_EL PELIGRO",105,"LP1912"),
    @("07:19:37","08:43","14_ABASTO",84,"LP1912"),
    @("06:57:11","08:54","17_ROMERO",117,"LP1912"),
    @("07:19:37","09:01","215A_EL PATO",102,"LP1912"),
    @("07:45:49","09:03","11_ETCHEVERRY",78,"LP1912"),
    @("07:19:37","09:10","16_P MOR-SANTA ANA",111,"LP1912"),
    @("07:19:37","09:16","27_EL RETIRO",117,"LP1912"),
    @("07:45:49","09:21","26_HERNANDEZ",96,"LP1912"),
    @("07:45:49","09:22","17_ROMERO",97,"LP1912"),
    @("07:45:49","09:23","11_ETCHEVERRY",98,"LP1912"),
    @("07:45:49","09:32","15_ABASTO",107,"LP1912"),
    @("07:45:49","09:33","10_OLMOS",108,"LP1912"),
    @("07:45:49","09:41","215C_EL PATO",116,"LP1912")
)
$ws.Range("A2").Value = "Última actualización: 07:45:49"
$ws.Range("A3").Value = "Total filas: 81"
for ($i = 0; $i -lt $data1.Count; $i++) {
    $r = $i + 6
    $row = $data1[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$data2 = @(
    @("03:55:13","04:34","215_ALUAR",39,"LP1912"),
    @("04:18:52","04:45","215A_EL PATO",27,"LP1912"),
    @("04:18:52","05:34","215B_EL PATO",76,"LP1912"),
    @("03:55:13","05:35","215B_EL PATO",100,"LP1912"),
    @("04:18:52","06:11","215A_EL PATO",113,"LP1912"),
    @("04:54:17","06:46","215C_EL PATO",112,"LP1912"),
    @("05:49:10","06:47","215C_EL PATO",58,"LP1912"),
    @("05:19:24","07:11","215A_EL PATO",112,"LP1912"),
    @("06:14:19","07:51","215D_EL PATO",97,"LP1912"),
    @("06:43:12","07:52","215D_EL PATO",69,"LP1912"),
    @("06:43:12","08:23","215B_EL PATO",100,"LP1912"),
    @("07:19:37","09:01","215A_EL PATO",102,"LP1912"),
    @("07:45:49","09:41","215C_EL PATO",116,"LP1912")
)
$ws.Range("A2").Value = "Última actualización: 07:45:49"
$ws.Range("A3").Value = "Total filas: 13"
for ($i = 0; $i -lt $data2.Count; $i++) {
    $r = $i + 6
    $row = $data2[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$data3 = @(
    @("04:18:52","05:43","215A_LA PLATA",85,"L6173"),
    @("03:55:13","05:44","215A_LA PLATA",109,"L6173"),
    @("04:18:52","06:08","215A_LA PLATA",110,"L6173"),
    @("04:54:17","06:09","215A_LA PLATA",75,"L6173"),
    @("04:40:48","06:32","215C_LA PLATA",112,"L6203"),
    @("04:54:17","06:33","215C_LA PLATA",99,"L6203"),
    @("05:19:24","07:00","215B_LP-P MOR-1 Y 57",101,"L6173"),
    @("05:49:10","07:35","215A_LA PLATA",106,"L6173"),
    @("07:19:37","07:38","215A_LA PLATA",19,"L6173"),
    @("06:14:19","08:07","215C_LA PLATA",113,"L6203"),
    @("07:19:37","08:08","215C_LA PLATA",49,"L6203"),
    @("07:45:49","08:09","215C_LA PLATA",24,"L6203"),
    @("06:57:11","08:10","215C_LA PLATA",73,"L6203"),
    @("07:19:37","08:35","215A_LA PLATA",76,"L6173"),
    @("06:57:11","08:38","215A_LA PLATA",101,"L6173"),
    @("07:19:37","09:08","215D_LA PLATA",109,"L6203")
)
$ws.Range("A2").Value = "Última actualización: 07:45:49"
$ws.Range("A3").Value = "Total filas: 16"
for ($i = 0; $i -lt $data3.Count; $i++) {
    $r = $i + 6
    $row = $data3[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
